# Edit script: update title, bullet lists, and meta description text
# for the Mystic Fortune Deluxe review document.
#
# We use Range.InsertXML with a fully-formed WordprocessingML package
# fragment (rather than Find/Replace) so that each paragraph's existing
# run layout -- including any leading empty <w:r/> run used for cursor
# placement / formatting marks -- is reproduced exactly, with only the
# visible text content changed.

$d = $word.ActiveDocument

function New-WordXmlPackage([string]$bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphXml($paragraph, [string]$innerParagraphXml) {
    $pkg = New-WordXmlPackage $innerParagraphXml
    $paragraph.Range.InsertXML($pkg) | Out-Null
}

# --- 1. Headline (Heading1, paragraph #1): drop the " - Slot Game Review" suffix
Set-ParagraphXml $d.Paragraphs(1) '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Mystic Fortune Deluxe for Free</w:t></w:r></w:p>'

# --- "What we like" bullets (paragraphs #34-#37)
Set-ParagraphXml $d.Paragraphs(34) '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Beautiful Chinese-inspired visuals</w:t></w:r></w:p>'
Set-ParagraphXml $d.Paragraphs(35) '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Extended Wild symbol for more winning combinations</w:t></w:r></w:p>'
Set-ParagraphXml $d.Paragraphs(36) '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Free spins and multipliers for increased rewards</w:t></w:r></w:p>'
Set-ParagraphXml $d.Paragraphs(37) '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Medium volatility for a balanced gameplay experience</w:t></w:r></w:p>'

# --- "What we don't like" bullets (paragraphs #39-#40)
Set-ParagraphXml $d.Paragraphs(39) '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Limited betting options</w:t></w:r></w:p>'
Set-ParagraphXml $d.Paragraphs(40) '<w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Moderate payout frequency</w:t></w:r></w:p>'

# --- Bold title repeated near the end (paragraph #41)
Set-ParagraphXml $d.Paragraphs(41) '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Mystic Fortune Deluxe for Free</w:t></w:r></w:p>'

# --- Italic meta description (paragraph #42)
Set-ParagraphXml $d.Paragraphs(42) '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the beautiful Chinese-inspired visuals and enjoy free spins in Mystic Fortune Deluxe.</w:t></w:r></w:p>'
